$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Row 3: update the battle-start message to include the set-bg tag inline.
$ws.Range("C3").Value = "[set-bg=bg_battle_2]Battle!"

# 2) Remove the old row 4 (separate DIALOGUE row that used to carry the
#    [set-bg=bg_battle_2]_ / _ filler line) - its content is now folded
#    into row 3 above, so the whole row goes away and everything below
#    shifts up by one.
$ws.Rows("4:4").Delete()

# 3) Insert the new battle-interrupt block (INTERRUPT / DIALOGUE / END_INTERRUPT)
#    right after the ENEMY list (which is now rows 9-12 post-shift).
$ws.Rows("13:15").Insert()

$ws.Range("A13").Value = "INTERRUPT"
$ws.Range("A15").Value = "END_INTERRUPT"

$ws.Range("B13").Value = "check-health,0,0.5"

$ws.Range("A14").Value = "DIALOGUE"
$ws.Range("B14").Value = "Tanuki"
$ws.Range("C14").Value = "Ow! That hurts!"
$ws.Range("D14").Value = "_"
$ws.Range("E14").Value = "tanuki"
$ws.Range("F14").Value = "CENTER"
$ws.Range("G14").Value = "END_DIALOGUE"

# Match the saved cursor/selection position from the authored workbook.
$ws.Range("C14").Select()
